$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.2380952380952381
$ws.Cells.Item(2, 3).Value = 0.4708994708994709
$ws.Cells.Item(2, 10).Value = 0.04232804232804233
$ws.Cells.Item(2, 16).Value = 0.1507936507936508
$ws.Cells.Item(2, 19).Value = 0.09788359788359788
$ws.Cells.Item(3, 2).Value = 0.0111731843575419
$ws.Cells.Item(3, 3).Value = 0.0223463687150838
$ws.Cells.Item(3, 10).Value = 0.07262569832402235
$ws.Cells.Item(3, 16).Value = 0.664804469273743
$ws.Cells.Item(3, 19).Value = 0.2290502793296089
$ws.Cells.Item(4, 10).Value = 0.04411764705882353
$ws.Cells.Item(4, 16).Value = 0.7205882352941176
$ws.Cells.Item(4, 19).Value = 0.2352941176470588
$ws.Cells.Item(5, 16).Value = 0.5
$ws.Cells.Item(5, 19).Value = 0.5
$ws.Cells.Item(6, 2).Value = 0.06598984771573604
$ws.Cells.Item(6, 4).Value = 0.01015228426395939
$ws.Cells.Item(6, 5).Value = 0.005076142131979695
$ws.Cells.Item(6, 6).Value = 0.04568527918781726
$ws.Cells.Item(6, 10).Value = 0.2131979695431472
$ws.Cells.Item(6, 15).Value = 0.03553299492385787
$ws.Cells.Item(6, 17).Value = 0.1878172588832487
$ws.Cells.Item(6, 18).Value = 0.05583756345177665
$ws.Cells.Item(6, 19).Value = 0.3807106598984771
$ws.Cells.Item(7, 2).Value = 0.1262626262626263
$ws.Cells.Item(7, 4).Value = 0.03535353535353535
$ws.Cells.Item(7, 6).Value = 0.0505050505050505
$ws.Cells.Item(7, 10).Value = 0.1414141414141414
$ws.Cells.Item(7, 15).Value = 0.005050505050505051
$ws.Cells.Item(7, 17).Value = 0.1767676767676768
$ws.Cells.Item(7, 18).Value = 0.04040404040404041
$ws.Cells.Item(7, 19).Value = 0.4242424242424243
$ws.Cells.Item(8, 2).Value = 0.1024590163934426
$ws.Cells.Item(8, 4).Value = 0.04098360655737705
$ws.Cells.Item(8, 6).Value = 0.06762295081967214
$ws.Cells.Item(8, 10).Value = 0.09016393442622951
$ws.Cells.Item(8, 15).Value = 0.006147540983606557
$ws.Cells.Item(8, 17).Value = 0.1905737704918033
$ws.Cells.Item(8, 18).Value = 0.1127049180327869
$ws.Cells.Item(8, 19).Value = 0.389344262295082
$ws.Cells.Item(9, 2).Value = 0.09727626459143969
$ws.Cells.Item(9, 4).Value = 0.02723735408560311
$ws.Cells.Item(9, 5).Value = 0.003891050583657588
$ws.Cells.Item(9, 6).Value = 0.0622568093385214
$ws.Cells.Item(9, 10).Value = 0.1206225680933852
$ws.Cells.Item(9, 15).Value = 0.007782101167315175
$ws.Cells.Item(9, 17).Value = 0.198443579766537
$ws.Cells.Item(9, 18).Value = 0.0933852140077821
$ws.Cells.Item(9, 19).Value = 0.3891050583657588
$ws.Cells.Item(10, 2).Value = 0.1319444444444444
$ws.Cells.Item(10, 4).Value = 0.02469135802469136
$ws.Cells.Item(10, 6).Value = 0.06327160493827161
$ws.Cells.Item(10, 10).Value = 0.1126543209876543
$ws.Cells.Item(10, 15).Value = 0.02777777777777778
$ws.Cells.Item(10, 17).Value = 0.2137345679012346
$ws.Cells.Item(10, 18).Value = 0.07407407407407407
$ws.Cells.Item(10, 19).Value = 0.3518518518518519
$ws.Cells.Item(11, 7).Value = 0.1608832807570978
$ws.Cells.Item(11, 10).Value = 0.09779179810725552
$ws.Cells.Item(11, 11).Value = 0.2271293375394322
$ws.Cells.Item(11, 12).Value = 0.4921135646687697
$ws.Cells.Item(11, 19).Value = 0.0220820189274448
$ws.Cells.Item(12, 7).Value = 0.7025316455696202
$ws.Cells.Item(12, 10).Value = 0.2088607594936709
$ws.Cells.Item(12, 12).Value = 0.0189873417721519
$ws.Cells.Item(12, 19).Value = 0.06962025316455696
$ws.Cells.Item(13, 7).Value = 0.6964285714285714
$ws.Cells.Item(13, 10).Value = 0.2678571428571428
$ws.Cells.Item(13, 19).Value = 0.03571428571428571
$ws.Cells.Item(15, 6).Value = 0.004347826086956522
$ws.Cells.Item(15, 8).Value = 0.1391304347826087
$ws.Cells.Item(15, 9).Value = 0.0782608695652174
$ws.Cells.Item(15, 10).Value = 0.391304347826087
$ws.Cells.Item(15, 11).Value = 0.05217391304347826
$ws.Cells.Item(15, 13).Value = 0.01304347826086956
$ws.Cells.Item(15, 15).Value = 0.04782608695652174
$ws.Cells.Item(15, 19).Value = 0.2739130434782608
$ws.Cells.Item(16, 6).Value = 0.01327433628318584
$ws.Cells.Item(16, 8).Value = 0.1725663716814159
$ws.Cells.Item(16, 9).Value = 0.05752212389380531
$ws.Cells.Item(16, 10).Value = 0.4026548672566372
$ws.Cells.Item(16, 11).Value = 0.163716814159292
$ws.Cells.Item(16, 13).Value = 0.02212389380530973
$ws.Cells.Item(16, 15).Value = 0.05309734513274336
$ws.Cells.Item(16, 19).Value = 0.1150442477876106
$ws.Cells.Item(17, 6).Value = 0.01020408163265306
$ws.Cells.Item(17, 8).Value = 0.1693877551020408
$ws.Cells.Item(17, 9).Value = 0.1306122448979592
$ws.Cells.Item(17, 10).Value = 0.4081632653061225
$ws.Cells.Item(17, 11).Value = 0.08571428571428572
$ws.Cells.Item(17, 13).Value = 0.01428571428571429
$ws.Cells.Item(17, 14).Value = 0.002040816326530612
$ws.Cells.Item(17, 15).Value = 0.06326530612244897
$ws.Cells.Item(17, 19).Value = 0.1163265306122449
$ws.Cells.Item(18, 6).Value = 0.005128205128205128
$ws.Cells.Item(18, 8).Value = 0.1487179487179487
$ws.Cells.Item(18, 9).Value = 0.1230769230769231
$ws.Cells.Item(18, 10).Value = 0.3487179487179487
$ws.Cells.Item(18, 11).Value = 0.1333333333333333
$ws.Cells.Item(18, 13).Value = 0.03589743589743589
$ws.Cells.Item(18, 15).Value = 0.05128205128205128
$ws.Cells.Item(18, 19).Value = 0.1538461538461539
$ws.Cells.Item(19, 6).Value = 0.009015777610818933
$ws.Cells.Item(19, 8).Value = 0.2291510142749812
$ws.Cells.Item(19, 9).Value = 0.1051840721262209
$ws.Cells.Item(19, 10).Value = 0.3516153268219384
$ws.Cells.Item(19, 11).Value = 0.09767092411720511
$ws.Cells.Item(19, 13).Value = 0.0270473328324568
$ws.Cells.Item(19, 14).Value = 0.0007513148009015778
$ws.Cells.Item(19, 15).Value = 0.06611570247933884
$ws.Cells.Item(19, 19).Value = 0.1134485349361382